$wb = $excel.ActiveWorkbook

$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.356209313247299
$wsLinear.Range("B3").Value = -0.01881335553529419
$wsLinear.Range("B4").Value = 5093.173399447834
$wsLinear.Range("B5").Value = "[1.0, 0.1674498097298499, -0.013183671409045911, -0.04628084952828638, -0.029934322614901202, -0.016852354190925668, 0.1729298857312265, 0.3035490714844914, 0.18262593962261844, 0.005472362640659067, -0.05167763796870478, -0.031456656675643464, -0.018616219420696006, 0.17558604414281315, 0.3479117965367412, 0.13423120701662838, -0.031014204624913354, -0.048288995162800934, -0.0474833872977779, -0.036178465186397106]"

$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.8265605979962498
$wsNonLinear.Range("B5").Value = -0.08517015475193113
$wsNonLinear.Range("B6").Value = 5164.906242224657
$wsNonLinear.Range("B7").Value = -1.996784302007529
$wsNonLinear.Range("B8").Value = 0.07353342375068327
$wsNonLinear.Range("B9").Value = 5020.679346452037
$wsNonLinear.Range("B10").Value = "[1.0, 0.16811469737224494, -0.008937768644628155, -0.04090467244627579, -0.025171643756092973, -0.01303701346772409, 0.17245055597509754, 0.29912170194702986, 0.18196495549946684, 0.00730175862572339, -0.04729621091876375, -0.02539587997651757, -0.016898418306406933, 0.1754246221566118, 0.34318480086001524, 0.1317063186275602, -0.02839029055072811, -0.04475012987006419, -0.04380919392537788, -0.03204470106182204]"
